# Auto-update stock values: 2025-12-12 07:55:10 UTC
$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ idx = 2; newCol = 73; oldCol = 72; isText = $false; header = 20251212; data = @(15510,57435,18170,15960,17125,19730,18520,1581,17215,5620,6120,6200,1761,19505,11355,6740,18295,16285) },
    @{ idx = 3; newCol = 73; oldCol = 72; isText = $false; header = 20251212; data = @(15930,57640,18315,15960,17125,19855,18655,1615,17285,5680,6180,6250,1774,19835,11465,6795,18295,16440) },
    @{ idx = 4; newCol = 73; oldCol = 72; isText = $false; header = 20251212; data = @(15405,57000,17785,15690,16855,19555,18415,1563,16747,5570,6035,6090,1732,19140,11280,6660,17765,15935) },
    @{ idx = 5; newCol = 73; oldCol = 72; isText = $false; header = 20251212; data = @(15470,57375,17810,15905,16985,19795,18625,1608,16885,5645,6135,6200,1767,19265,11430,6770,17845,15965) },
    @{ idx = 6; newCol = 73; oldCol = 72; isText = $false; header = 20251212; data = @(118334,890897,214682,167695,1743091,1190061,2682396,3121846,858919,374862,1853968,5108733,26813627,1379949,304760,33097,90188,45934) },
    @{ idx = 7; newCol = 54; oldCol = 53; isText = $false; header = 20251212; data = @(31,79,0,50,83,95,87,75,54,86,87,88,87,40,88,85,54,50) },
    @{ idx = 8; newCol = 14; oldCol = 13; isText = $false; header = 20251212; data = @(72,80,61,69,86,93,86,61,70,93,93,94,70,81,86,82,83,81) },
    @{ idx = 9; newCol = 54; oldCol = 53; isText = $false; header = 20251212; data = @(-49,33,-78,-6,32,56,45,43,-23,61,63,65,51,-14,62,56,-27,-13) },
    @{ idx = 10; newCol = 14; oldCol = 13; isText = $false; header = 20251212; data = @(41,34,21,14,42,53,40,27,19,55,56,62,36,58,48,44,38,58) },
    @{ idx = 11; newCol = 54; oldCol = 53; isText = $true; header = 20251212; data = @(97,102,96,100,102,104,103,105,99,107,108,110,107,99,106,106,98,99) },
    @{ idx = 12; newCol = 35; oldCol = 34; isText = $true; header = 20251212; data = @(-40.21,-16.48,-34.52,-22.99,-7.11,5.44,-12.85,-31.82,-7.66,29.61,30.31,45.15,-19.86,-52.09,5.77,0.71,-9.88,-48.54) },
    @{ idx = 13; newCol = 14; oldCol = 13; isText = $true; header = 20251212; data = @(53,19,68,10,43,15,19,22,29,15,33,38,26,42,25,7,25,48) }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.idx)

    # Copy the last existing header cell's formatting onto the new header
    # cell so the header style ("s" attribute) matches the rest of the row.
    $ws.Cells.Item(1, $s.oldCol).Copy()
    $ws.Cells.Item(1, $s.newCol).PasteSpecial(-4122)

    # Match the new column's width to its left neighbour.
    $srcWidth = $ws.Columns.Item($s.oldCol).ColumnWidth
    $ws.Columns.Item($s.newCol).ColumnWidth = $srcWidth

    # Header cell (row 1): the date value. Some sheets store the header row
    # as text (matches existing convention on that sheet), others as numbers.
    $headerCell = $ws.Cells.Item(1, $s.newCol)
    if ($s.isText) {
        $headerCell.Value = "'" + $s.header
    } else {
        $headerCell.Value = $s.header
    }

    # Data rows 2-19
    $r = 2
    foreach ($v in $s.data) {
        $ws.Cells.Item($r, $s.newCol).Value = $v
        $r = $r + 1
    }
}

